# Update the organisational practices model fit results (Sheet1, column B)
# to reflect a re-fit of the model (npar changed 6 -> 7 and all downstream
# statistics updated accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = 7.0
$ws.Range("B3").Value  = [double]"1.2710865025773388E-16"
$ws.Range("B4").Value  = [double]"3.894863261197482E-12"
$ws.Range("B7").Value  = [double]"3.894863261197482E-12"
$ws.Range("B11").Value = 17921.096717470624
$ws.Range("B14").Value = 17922.266501851664
$ws.Range("B17").Value = 0.9999347301090005
$ws.Range("B18").Value = 0.9999999999999998
$ws.Range("B25").Value = 0.9999999999999998
$ws.Range("B26").Value = 0.9999999999999998
$ws.Range("B35").Value = 0.9999999999999998
$ws.Range("B49").Value = [double]"4.9452456701650085E-9"
$ws.Range("B50").Value = [double]"6.056664272306056E-9"
$ws.Range("B51").Value = [double]"6.056664272306056E-9"
$ws.Range("B52").Value = [double]"4.9452456701650085E-9"
$ws.Range("B53").Value = [double]"6.056664272306056E-9"
$ws.Range("B54").Value = [double]"6.056664272306056E-9"
$ws.Range("B55").Value = [double]"8.565416756635797E-9"
$ws.Range("B60").Value = 0.9999999999999998
